$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("preseaon_review")

$ws.Range("A4").Value = 'Huang2020'
$ws.Range("B4").Value = 'Phenological changes in herbaceous plants in ChinaÕs grasslands and their responses to climate change: a meta-analysis'
$ws.Range("C4").Value = 'N'
$ws.Range("D4").Value = 'metaanalysis with multiple preseasons in different studies'

$ws.Range("A5").Value = 'Guo2020'
$ws.Range("B5").Value = 'Impact of spring phenology variation on GPP and its lag feedback for winter wheat over the North China Plain'
$ws.Range("C5").Value = 'Y'
$ws.Range("E5").Value = 'parcorr'
$ws.Range("F5").Value = '0-3 month'
$ws.Range("G5").Value = 'multiyear average Green up date'
$ws.Range("I5").Value = '"The candidate previous month with the largest absolute partial correlation coefficient was finally selected as the preseason."'

$ws.Range("A6").Value = 'Yu2020'
$ws.Range("B6").Value = 'Forest Phenology Shifts in Response to Climate Change over China-Mongolia-Russia International Economic Corridor'
$ws.Range("C6").Value = 'N'
$ws.Range("D6").Value = 'preseason Snow Water Equivalent (SWE'

$ws.Range("A7").Value = 'Zhou2020'
$ws.Range("B7").Value = 'Assessment of varying changes of vegetation and the response to climatic factors using GIMMS NDVI3g on the Tibetan Plateau'
$ws.Range("C7").Value = 'Y'
$ws.Range("F7").Value = 'Jan'
$ws.Range("G7").Value = 'April'

$ws.Range("A8").Value = 'Huang2020b'
$ws.Range("B8").Value = 'Effect of preseason diurnal temperature range on the start of vegetation growing season in the Northern Hemisphere'
$ws.Range("C8").Value = 'Y'
$ws.Range("E8").Value = 'parcorr'
$ws.Range("F8").Value = 'Jan'
$ws.Range("G8").Value = 'Month of Start of Spring'
$ws.Range("I8").Value = 'Candidate preseason periods were evaluated in one-month intervals working backward in time from the month of SOS to January of the current year. A partial coefficient between DTR and SOS was calculated for each candidate period, and the period with the largest partial correlation coefficient (absolute value) was selected as the preseason period. The month of SOS was determined from the multiyear average of SOS dates. If the average date was in the second half of the month, that month was chosen as the SOS month, and if the average date was in the first half of the month, the previous month was chosen.'

$ws.Range("A9").Value = 'Chai2020'
$ws.Range("B9").Value = 'The relative controls of temperature and soil moisture on the start of carbon flux phenology and net ecosystem production in two alpine meadows on the Qinghai-Tibetan Plateau'

$ws.Range("A10").Value = 'An2020'
$ws.Range("B10").Value = 'Precipitation and Minimum Temperature are Primary Climatic Controls of Alpine Grassland Autumn Phenology on the Qinghai-Tibet Plateau'
$ws.Range("C10").Value = 'N'
$ws.Range("D10").Value = 'autumn phenology'

$ws.Range("A11").Value = 'Li2020'
$ws.Range("B11").Value = 'Change in Autumn Vegetation Phenology and the Climate Controls From 1982 to 2012 on the Qinghai-Tibet Plateau'
$ws.Range("C11").Value = 'N'
$ws.Range("D11").Value = 'autumn phenology'

# Activate the preseaon_review sheet/tab and set the selection to D12
$ws.Select()
$ws.Range("D12").Select()
